$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.2957516057923159
$ws.Range("D2").Value = 0.7701933032437294

$ws.Range("C3").Value = -0.4501836421180653
$ws.Range("D3").Value = 0.6569828885183546

$ws.Range("C4").Value = 1.76383606711743
$ws.Range("D4").Value = 0.09164360164728524

$ws.Range("C5").Value = 0.4769451174299562
$ws.Range("D5").Value = 0.6381047597614056

$ws.Range("C6").Value = -0.087852129854519
$ws.Range("D6").Value = 0.9307889843913399

$ws.Range("C7").Value = 2.251226753375598
$ws.Range("D7").Value = 0.03469407328537821

$ws.Range("C8").Value = 0.6479829579037554
$ws.Range("D8").Value = 0.5237017288343173

$ws.Range("C9").Value = 2.52898627634921
$ws.Range("D9").Value = 0.01911523949051164

$ws.Range("C10").Value = 1.224049879130315
$ws.Range("D10").Value = 0.2338809143816261

$ws.Range("C11").Value = -1.539566243516282
$ws.Range("D11").Value = 0.1379268892941832
